$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '28.754.22'
$ws.Cells.Item(2, 5).Value = '  -2.14%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '1.853.21'
$ws.Cells.Item(3, 5).Value = '  -3.15%  '
# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.009'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.34%  '
# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '337.01'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.90%  '
# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.05%  '
# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4679'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -2.37%  '
# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3927'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.82%  '
# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '46.94'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.89%  '
# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07921'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.05%  '
# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9802'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.10%  '
# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.39'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -3.85%  '
# Row 13
$ws.Cells.Item(13, 4).Value = '1.837.48'
$ws.Cells.Item(13, 5).Value = '  -4.29%  '
# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.840'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -3.04%  '
# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.981'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -3.01%  '
# Row 16
$ws.Cells.Item(16, 2).Value = 'TRON'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06805'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.39%  '
# Row 17
$ws.Cells.Item(17, 2).Value = 'BinanceUSD'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.011'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.15%  '
# Row 18
$ws.Cells.Item(18, 2).Value = 'Litecoin'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '87.83'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -3.72%  '
# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00001009'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -2.67%  '
# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.03'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.64%  '
# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.02%  '
# Row 22
$ws.Cells.Item(22, 4).Value = '28.792.85'
$ws.Cells.Item(22, 5).Value = '  -2.04%  '
# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.436'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -3.72%  '
# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.36'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -3.45%  '
# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.138'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.44%  '
# Row 26
$ws.Cells.Item(26, 4).Value = '2.125.67'
$ws.Cells.Item(26, 5).Value = '  -1.09%  '
# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '154.28'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.85%  '
# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.362'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -3.65%  '
# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.46'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -2.46%  '
# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.016'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -3.65%  '
# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '117.36'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -2.33%  '
# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9757'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.22%  '
# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.09443'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.69%  '
# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.406'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.38%  '
# Row 35
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.510'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.28%  '
# Row 36
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.359'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.35%  '
# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06193'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -3.07%  '
# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.02203'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -3.16%  '
# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.163'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.61%  '
# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5727'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -3.11%  '
# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.644'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -3.34%  '
# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.21'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.54%  '
# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1792'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.63%  '
# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.433'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.08%  '
# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.233'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -3.68%  '
# Row 46
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5409'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.49%  '
# Row 47
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.84'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.44%  '
# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07175'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -3.86%  '
# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.917'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.23%  '
# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '115.54'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.79%  '
# Row 51
$ws.Cells.Item(51, 2).Value = 'PaxDollar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.011'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.09%  '
